$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new data rows before the existing row 411 (a new weekly price report
# block for Fruta/Manzana at Terminal Hortofrutícola Agro Chillán), pushing the
# rows that were 411-443 down to 415-447.
$ws.Rows("411:414").Insert()

# Common / constant fields shared by the whole table block (same as neighboring rows)
$commonA = 7
$commonB = "Terminal Hortofrutícola Agro Chillán"
$commonC = "Ñuble"
$commonE = 16
$commonF = "Fruta"
$commonG = 100104
$commonH = "Frutos de pepita"
$commonI = 100104002
$commonJ = "Manzana"
$commonQ = "$/caja 16 kilos empedrada"
$commonR = "Provincia de Curicó"
$commonT = 16

# --- Row 411: Fuji royal / Primera ---
$ws.Range("A411").Value = $commonA
$ws.Range("B411").Value = $commonB
$ws.Range("C411").Value = $commonC
$ws.Range("D411").Value = 44461
$ws.Range("E411").Value = $commonE
$ws.Range("F411").Value = $commonF
$ws.Range("G411").Value = $commonG
$ws.Range("H411").Value = $commonH
$ws.Range("I411").Value = $commonI
$ws.Range("J411").Value = $commonJ
$ws.Range("K411").Value = "Fuji royal"
$ws.Range("L411").Value = "Primera"
$ws.Range("M411").Value = 240
$ws.Range("N411").Value = 8500
$ws.Range("O411").Value = 9000
$ws.Range("P411").Value = 8750
$ws.Range("Q411").Value = $commonQ
$ws.Range("R411").Value = $commonR
$ws.Range("S411").Value = 547
$ws.Range("T411").Value = $commonT

# --- Row 412: Fuji royal / Segunda ---
$ws.Range("A412").Value = $commonA
$ws.Range("B412").Value = $commonB
$ws.Range("C412").Value = $commonC
$ws.Range("D412").Value = 44461
$ws.Range("E412").Value = $commonE
$ws.Range("F412").Value = $commonF
$ws.Range("G412").Value = $commonG
$ws.Range("H412").Value = $commonH
$ws.Range("I412").Value = $commonI
$ws.Range("J412").Value = $commonJ
$ws.Range("K412").Value = "Fuji royal"
$ws.Range("L412").Value = "Segunda"
$ws.Range("M412").Value = 90
$ws.Range("N412").Value = 8000
$ws.Range("O412").Value = 8000
$ws.Range("P412").Value = 8000
$ws.Range("Q412").Value = $commonQ
$ws.Range("R412").Value = $commonR
$ws.Range("S412").Value = 500
$ws.Range("T412").Value = $commonT

# --- Row 413: Granny Smith / Primera ---
$ws.Range("A413").Value = $commonA
$ws.Range("B413").Value = $commonB
$ws.Range("C413").Value = $commonC
$ws.Range("D413").Value = 44461
$ws.Range("E413").Value = $commonE
$ws.Range("F413").Value = $commonF
$ws.Range("G413").Value = $commonG
$ws.Range("H413").Value = $commonH
$ws.Range("I413").Value = $commonI
$ws.Range("J413").Value = $commonJ
$ws.Range("K413").Value = "Granny Smith"
$ws.Range("L413").Value = "Primera"
$ws.Range("M413").Value = 240
$ws.Range("N413").Value = 8500
$ws.Range("O413").Value = 9000
$ws.Range("P413").Value = 8750
$ws.Range("Q413").Value = $commonQ
$ws.Range("R413").Value = $commonR
$ws.Range("S413").Value = 547
$ws.Range("T413").Value = $commonT

# --- Row 414: Granny Smith / Segunda ---
$ws.Range("A414").Value = $commonA
$ws.Range("B414").Value = $commonB
$ws.Range("C414").Value = $commonC
$ws.Range("D414").Value = 44461
$ws.Range("E414").Value = $commonE
$ws.Range("F414").Value = $commonF
$ws.Range("G414").Value = $commonG
$ws.Range("H414").Value = $commonH
$ws.Range("I414").Value = $commonI
$ws.Range("J414").Value = $commonJ
$ws.Range("K414").Value = "Granny Smith"
$ws.Range("L414").Value = "Segunda"
$ws.Range("M414").Value = 180
$ws.Range("N414").Value = 7500
$ws.Range("O414").Value = 8000
$ws.Range("P414").Value = 7750
$ws.Range("Q414").Value = $commonQ
$ws.Range("R414").Value = $commonR
$ws.Range("S414").Value = 484
$ws.Range("T414").Value = $commonT
